# This workbook tracks weekly price observations for garlic ("Ajo") at
# "Vega Monumental Concepción". A new weekly record is inserted at the
# top of the data block (row 142), pushing all existing rows down by
# one. The new record reuses most attributes of the prior top row but
# carries a newer date and updated volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 142; this shifts rows 142:168
# down to 143:169 (and all their formatting/values move with them),
# matching the dimension growing from R168 to R169.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new weekly observation.
$ws.Range("A142").Value = 11
$ws.Range("B142").Value = "Vega Monumental Concepción"
$ws.Range("C142").Value = "Bíobío"
$ws.Range("D142").Value = 44694
$ws.Range("E142").Value = 8
$ws.Range("F142").Value = 100112003
$ws.Range("G142").Value = "Ajo"
$ws.Range("H142").Value = "Chino"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 400
$ws.Range("K142").Value = 18000
$ws.Range("L142").Value = 19000
$ws.Range("M142").Value = 18500
$ws.Range("N142").Value = "$/caja 10 kilos"
$ws.Range("O142").Value = "China"
$ws.Range("P142").Value = 1850
$ws.Range("Q142").Value = 10
$ws.Range("R142").Value = "Hortaliza"
